$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G column timestamps for existing rows 2-21
$ws.Cells.Item(2, 7).Value = "'2024-08-15 21:51:18"
$ws.Cells.Item(3, 7).Value = "'2024-08-15 21:51:19"
$ws.Cells.Item(4, 7).Value = "'2024-08-15 21:51:20"
$ws.Cells.Item(5, 7).Value = "'2024-08-15 21:51:22"
$ws.Cells.Item(6, 7).Value = "'2024-08-15 21:51:23"
$ws.Cells.Item(7, 7).Value = "'2024-08-15 21:51:24"
$ws.Cells.Item(8, 7).Value = "'2024-08-15 21:51:25"
$ws.Cells.Item(9, 7).Value = "'2024-08-15 21:51:26"
$ws.Cells.Item(10, 7).Value = "'2024-08-15 21:51:27"
$ws.Cells.Item(11, 7).Value = "'2024-08-15 21:51:28"
$ws.Cells.Item(12, 7).Value = "'2024-08-15 21:51:29"
$ws.Cells.Item(13, 7).Value = "'2024-08-15 21:51:30"
$ws.Cells.Item(14, 7).Value = "'2024-08-15 21:51:31"
$ws.Cells.Item(15, 7).Value = "'2024-08-15 21:51:32"
$ws.Cells.Item(16, 7).Value = "'2024-08-15 21:51:33"
$ws.Cells.Item(17, 7).Value = "'2024-08-15 21:51:35"
$ws.Cells.Item(18, 7).Value = "'2024-08-15 21:51:36"
$ws.Cells.Item(19, 7).Value = "'2024-08-15 21:51:37"
$ws.Cells.Item(20, 7).Value = "'2024-08-15 21:51:38"
$ws.Cells.Item(21, 7).Value = "'2024-08-15 21:51:39"

# Append new rows 22-41
$newRows = @(
    @{ r=22; A="Appartement à Genève"; B="CHF 1,430.– / month"; C="2"; D="25m²"; E="Rue de la Navigation, 1201 Genève"; F="https://www.homegate.ch/rent/4001355523"; G="2024-08-15 21:52:00" },
    @{ r=23; A="Studio à Genève - meublé - piscine - proche des organisations internationales"; B="CHF 1,760.– / month"; C="1"; D="26m²"; E="Avenue De-Budé, 1202 Genève"; F="https://www.homegate.ch/rent/4001352171"; G="2024-08-15 21:52:01" },
    @{ r=24; A="Appartement à Genève"; B="CHF 1,500.– / month"; C="1.5"; D="25m²"; E="Rue Louise-De-Frotté 35, 1205 Genève"; F="https://www.homegate.ch/rent/4001337045"; G="2024-08-15 21:52:02" },
    @{ r=25; A="City-Penthouse en plein coeur des rues basses"; B="CHF 4,850.– / month"; C="3"; D="100m²"; E="Rue de la Croix-d'Or 27, 1204 Genève"; F="https://www.homegate.ch/rent/4001335811"; G="2024-08-15 21:52:03" },
    @{ r=26; A="Appartement en campagne"; B="CHF 2,210.– / month"; C="4"; D="78m²"; E="Route de Mon-Idée 49, 1226 Thônex"; F="https://www.homegate.ch/rent/4001333875"; G="2024-08-15 21:52:04" },
    @{ r=27; A="Appartement à Thônex"; B="CHF 4,564.– / month"; C="6"; D="139m²"; E="Cour de l'Emine 1, 1226 Thônex"; F="https://www.homegate.ch/rent/4001328477"; G="2024-08-15 21:52:05" },
    @{ r=28; A="Splendide logement dans le prestigieux quartier de Champel"; B="CHF 2,900.– / month"; C="3"; D="69m²"; E="Rue de Contamines 35, 1206 Genève"; F="https://www.homegate.ch/rent/4001313855"; G="2024-08-15 21:52:06" },
    @{ r=29; A="5 pièces de haute qualité architecturale, centre ville"; B="CHF 4,250.– / month"; C="5"; D="100m²"; E="Avenue de Frontenex 4, 1207 Genève"; F="https://www.homegate.ch/rent/4001313160"; G="2024-08-15 21:52:07" },
    @{ r=30; A="Appartement de 5.5 pièces meublé à Thônex"; B="CHF 3,900.– / month"; C="5.5"; D="130m²"; E="Chemin Etienne-Chennaz 15, 1226 Thônex"; F="https://www.homegate.ch/rent/4001308653"; G="2024-08-15 21:52:08" },
    @{ r=31; A="Spacious, peaceful, nicely furnished Apartment close to the UN"; B="CHF 4,775.– / month"; C="4.5"; D="120m²"; E="Chemin du Point-du-Jour, 1202 Geneva"; F="https://www.homegate.ch/rent/4001305957"; G="2024-08-15 21:52:09" },
    @{ r=32; A="Splendide logement dans le prestigieux quartier de Champel"; B="CHF 2,850.– / month"; C="3"; D="80m²"; E="Rue de Contamines 35, 1206 Genève"; F="https://www.homegate.ch/rent/4001287129"; G="2024-08-15 21:52:10" },
    @{ r=33; A="Beau logement 4.5 pièces très spacieux proche de Cornavin / Mercier"; B="CHF 2,920.– / month"; C="4.5"; D="N/A"; E="Rue Jean-Gutenberg, 1201 Genève"; F="https://www.homegate.ch/rent/4001271361"; G="2024-08-15 21:52:12" },
    @{ r=34; A="Appartement à Genève"; B="CHF 1,270.– / month"; C="1.5"; D="30m²"; E="Rue Jean-Charles Amat 15, 1202 Genève"; F="https://www.homegate.ch/rent/4001243823"; G="2024-08-15 21:52:13" },
    @{ r=35; A="Splendide logement dans le quartier de Champel"; B="CHF 2,750.– / month"; C="3"; D="63m²"; E="Rue de Contamines 35, 1206 Genève"; F="https://www.homegate.ch/rent/4001232755"; G="2024-08-15 21:52:14" },
    @{ r=36; A="Studio à Chêne-Bougeries"; B="CHF 1,690.– / month"; C="1"; D="28m²"; E="Avenue Pierre-Odier, 1224 Chêne-Bougeries"; F="https://www.homegate.ch/rent/4001178602"; G="2024-08-15 21:52:15" },
    @{ r=37; A="Appartement à Genève"; B="CHF 1,600.– / month"; C="3"; D="60m²"; E="Rue de Contamines, 1206 Genève"; F="https://www.homegate.ch/rent/4001158284"; G="2024-08-15 21:52:16" },
    @{ r=38; A="Appartement de charme dans le vieux Chêne-Bourg"; B="CHF 2,090.– / month"; C="2.5"; D="56m²"; E="Rue du Gothard 7, 1225 Chêne-Bourg"; F="https://www.homegate.ch/rent/4000848882"; G="2024-08-15 21:52:17" },
    @{ r=39; A="Bel appartement lumineux"; B="CHF 2,255.– / month"; C="3"; D="65m²"; E="Rue des Bossons 26, 1213 Onex"; F="https://www.homegate.ch/rent/3003397551"; G="2024-08-15 21:52:18" },
    @{ r=40; A="Magnifique appartement au centre-ville de Genève ! Terme fixe au 31.08.2026"; B="CHF 4,870.– / month"; C="6"; D="105m²"; E="Rue Charles-Giron 14, 1203 Genève"; F="https://www.homegate.ch/rent/4001225099"; G="2024-08-15 21:52:19" },
    @{ r=41; A="Un projet unique au coeur de Genève"; B="CHF 5,460.– / month"; C="4.5"; D="104m²"; E="Rue du Vieux Collège 3, 1204 Genève"; F="https://www.homegate.ch/rent/4000656611"; G="2024-08-15 21:52:20" }
)

foreach ($row in $newRows) {
    $ws.Cells.Item($row.r, 1).Value = $row.A
    $ws.Cells.Item($row.r, 2).Value = $row.B
    $ws.Cells.Item($row.r, 3).Value = "'" + $row.C
    $ws.Cells.Item($row.r, 4).Value = $row.D
    $ws.Cells.Item($row.r, 5).Value = $row.E
    $ws.Cells.Item($row.r, 6).Value = $row.F
    $ws.Cells.Item($row.r, 7).Value = "'" + $row.G
}